# This workbook tracks FFXIV crafting-leve profitability per job (8 sheets:
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Columns H:N are market-board
# snapshot data (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# refreshed on a schedule from an external price API. This run refreshes the
# cached prices/profits for the rows the API returned new data for.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1232.5883  # H28: was 222.66667
$ws.Cells.Item(28, 9).Value = 227.23077  # I28: was 182.94737
$ws.Cells.Item(28, 10).Value = 4500  # J28: was 600
$ws.Cells.Item(28, 11).Value = 227.23077  # K28: was 182.94737
$ws.Cells.Item(28, 12).Value = 4500  # L28: was 600
$ws.Cells.Item(28, 13).Value = 257.76923  # M28: was 302.05263
$ws.Cells.Item(28, 14).Value = -5470  # N28: was -1570

$ws.Cells.Item(62, 8).Value = 8115.88  # H62: was 8412.375
$ws.Cells.Item(62, 9).Value = 7690  # I62: was 8108.125
$ws.Cells.Item(62, 11).Value = 7690  # K62: was 8108.125
$ws.Cells.Item(62, 13).Value = -7066  # M62: was -7484.125

$ws.Cells.Item(64, 8).Value = 3215.625  # H64: was 2978.5715
$ws.Cells.Item(64, 9).Value = 3516.6667  # I64: was 2927.2727
$ws.Cells.Item(64, 11).Value = 3516.6667  # K64: was 2927.2727
$ws.Cells.Item(64, 13).Value = -3268.6667  # M64: was -2679.2727

$ws.Cells.Item(65, 8).Value = 8115.88  # H65: was 8412.375
$ws.Cells.Item(65, 9).Value = 7690  # I65: was 8108.125
$ws.Cells.Item(65, 11).Value = 38450  # K65: was 40540.625
$ws.Cells.Item(65, 13).Value = -35330  # M65: was -37420.625

$ws.Cells.Item(67, 8).Value = 3215.625  # H67: was 2978.5715
$ws.Cells.Item(67, 9).Value = 3516.6667  # I67: was 2927.2727
$ws.Cells.Item(67, 11).Value = 3516.6667  # K67: was 2927.2727
$ws.Cells.Item(67, 13).Value = -2658.6667  # M67: was -2069.2727

$ws.Cells.Item(137, 8).Value = 1609.1724  # H137: was 1748.16
$ws.Cells.Item(137, 9).Value = 1582.421  # I137: was 1744
$ws.Cells.Item(137, 10).Value = 1660  # J137: was 1755.5555
$ws.Cells.Item(137, 11).Value = 4747.263  # K137: was 5232
$ws.Cells.Item(137, 12).Value = 4980  # L137: was 5266.666499999999
$ws.Cells.Item(137, 13).Value = -2197.263  # M137: was -2682
$ws.Cells.Item(137, 14).Value = -10080  # N137: was -10366.6665

$ws.Cells.Item(138, 8).Value = 2261.853  # H138: was 3139
$ws.Cells.Item(138, 9).Value = 1516.75  # I138: was 2953.2222
$ws.Cells.Item(138, 10).Value = 3326.2856  # J138: was 3222.6
$ws.Cells.Item(138, 11).Value = 4550.25  # K138: was 8859.6666
$ws.Cells.Item(138, 12).Value = 9978.856800000001  # L138: was 9667.799999999999
$ws.Cells.Item(138, 13).Value = 589.75  # M138: was -3719.6666
$ws.Cells.Item(138, 14).Value = -20258.8568  # N138: was -19947.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14709706  # H32: was 17548154
$ws.Cells.Item(32, 9).Value = 3734.5688  # I32: was 4188
$ws.Cells.Item(32, 10).Value = 100004344  # J32: was 125004936
$ws.Cells.Item(32, 11).Value = 3734.5688  # K32: was 4188
$ws.Cells.Item(32, 12).Value = 100004344  # L32: was 125004936
$ws.Cells.Item(32, 13).Value = -3447.5688  # M32: was -3901
$ws.Cells.Item(32, 14).Value = -100004918  # N32: was -125005510

$ws.Cells.Item(61, 8).Value = 4763041  # H61: was 6668033
$ws.Cells.Item(61, 9).Value = 5377521.5  # I61: was 6668033
$ws.Cells.Item(61, 10).Value = 820  # J61: was 0
$ws.Cells.Item(61, 11).Value = 5377521.5  # K61: was 6668033
$ws.Cells.Item(61, 12).Value = 820  # L61: was 0
$ws.Cells.Item(61, 13).Value = -5377309.5  # M61: was -6667821
$ws.Cells.Item(61, 14).Value = -1244  # N61: was empty

$ws.Cells.Item(74, 8).Value = 1017.67645  # H74: was 1367.0322
$ws.Cells.Item(74, 9).Value = 875.0833  # I74: was 1443.75
$ws.Cells.Item(74, 10).Value = 1359.9  # J74: was 1227.5454
$ws.Cells.Item(74, 11).Value = 875.0833  # K74: was 1443.75
$ws.Cells.Item(74, 12).Value = 1359.9  # L74: was 1227.5454
$ws.Cells.Item(74, 13).Value = -1.083300000000008  # M74: was -569.75
$ws.Cells.Item(74, 14).Value = -3107.9  # N74: was -2975.5454

$ws.Cells.Item(77, 8).Value = 1017.67645  # H77: was 1367.0322
$ws.Cells.Item(77, 9).Value = 875.0833  # I77: was 1443.75
$ws.Cells.Item(77, 10).Value = 1359.9  # J77: was 1227.5454
$ws.Cells.Item(77, 11).Value = 4375.4165  # K77: was 7218.75
$ws.Cells.Item(77, 12).Value = 6799.5  # L77: was 6137.727
$ws.Cells.Item(77, 13).Value = -7.416500000000269  # M77: was -2850.75
$ws.Cells.Item(77, 14).Value = -15535.5  # N77: was -14873.727

$ws.Cells.Item(88, 8).Value = 3571.4285  # H88: was 2789
$ws.Cells.Item(88, 9).Value = 4500  # I88: was 2959.1
$ws.Cells.Item(88, 10).Value = 2875  # J88: was 2546
$ws.Cells.Item(88, 11).Value = 4500  # K88: was 2959.1
$ws.Cells.Item(88, 12).Value = 2875  # L88: was 2546
$ws.Cells.Item(88, 13).Value = -4094  # M88: was -2553.1
$ws.Cells.Item(88, 14).Value = -3687  # N88: was -3358

$ws.Cells.Item(91, 8).Value = 3571.4285  # H91: was 2789
$ws.Cells.Item(91, 9).Value = 4500  # I91: was 2959.1
$ws.Cells.Item(91, 10).Value = 2875  # J91: was 2546
$ws.Cells.Item(91, 11).Value = 4500  # K91: was 2959.1
$ws.Cells.Item(91, 12).Value = 2875  # L91: was 2546
$ws.Cells.Item(91, 13).Value = -3096  # M91: was -1555.1
$ws.Cells.Item(91, 14).Value = -5683  # N91: was -5354

$ws.Cells.Item(92, 8).Value = 0  # H92: was 10550
$ws.Cells.Item(92, 10).Value = 0  # J92: was 10550
$ws.Cells.Item(92, 12).ClearContents()  # L92: was 10550 -> empty
$ws.Cells.Item(92, 14).Value = 0  # N92: was -15542

$ws.Cells.Item(132, 8).Value = 1091.4043  # H132: was 1540.0454
$ws.Cells.Item(132, 9).Value = 795.6111  # I132: was 992.6875
$ws.Cells.Item(132, 10).Value = 2059.4546  # J132: was 2999.6667
$ws.Cells.Item(132, 11).Value = 2386.8333  # K132: was 2978.0625
$ws.Cells.Item(132, 12).Value = 6178.3638  # L132: was 8999.000100000001
$ws.Cells.Item(132, 13).Value = 143.1667000000002  # M132: was -448.0625
$ws.Cells.Item(132, 14).Value = -11238.3638  # N132: was -14059.0001

$ws.Cells.Item(136, 8).Value = 4763041  # H136: was 6668033
$ws.Cells.Item(136, 9).Value = 5377521.5  # I136: was 6668033
$ws.Cells.Item(136, 10).Value = 820  # J136: was 0
$ws.Cells.Item(136, 11).Value = 16132564.5  # K136: was 20004099
$ws.Cells.Item(136, 12).Value = 2460  # L136: was 0
$ws.Cells.Item(136, 13).Value = -16130014.5  # M136: was -20001549
$ws.Cells.Item(136, 14).Value = -7560  # N136: was empty

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(23, 8).Value = 4996.6665  # H23: was 4933.3335
$ws.Cells.Item(23, 10).Value = 4996.6665  # J23: was 4933.3335
$ws.Cells.Item(23, 12).Value = 4996.6665  # L23: was 4933.3335
$ws.Cells.Item(23, 14).Value = -5562.6665  # N23: was -5499.3335

$ws.Cells.Item(99, 8).Value = 1786  # H99: was 2958.7144
$ws.Cells.Item(99, 9).Value = 808.9  # I99: was 1005
$ws.Cells.Item(99, 11).Value = 808.9  # K99: was 1005
$ws.Cells.Item(99, 13).Value = 689.1  # M99: was 493

$ws.Cells.Item(134, 8).Value = 2416462.8  # H134: was 2416467.5
$ws.Cells.Item(134, 9).Value = 912.6842  # I134: was 920.0513
$ws.Cells.Item(134, 10).Value = 13890326  # J134: was 15874518
$ws.Cells.Item(134, 11).Value = 2738.0526  # K134: was 2760.1539
$ws.Cells.Item(134, 12).Value = 41670978  # L134: was 47623554
$ws.Cells.Item(134, 13).Value = -203.0526  # M134: was -225.1538999999998
$ws.Cells.Item(134, 14).Value = -41676048  # N134: was -47628624

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1323613.8  # H31: was 1355890.9
$ws.Cells.Item(31, 9).Value = 1684303.5  # I31: was 1736934.2
$ws.Cells.Item(31, 10).Value = 1084.6666  # J31: was 1070.2222
$ws.Cells.Item(31, 11).Value = 1684303.5  # K31: was 1736934.2
$ws.Cells.Item(31, 12).Value = 1084.6666  # L31: was 1070.2222
$ws.Cells.Item(31, 13).Value = -1684008.5  # M31: was -1736639.2
$ws.Cells.Item(31, 14).Value = -1674.6666  # N31: was -1660.2222

$ws.Cells.Item(34, 8).Value = 1323613.8  # H34: was 1355890.9
$ws.Cells.Item(34, 9).Value = 1684303.5  # I34: was 1736934.2
$ws.Cells.Item(34, 10).Value = 1084.6666  # J34: was 1070.2222
$ws.Cells.Item(34, 11).Value = 1684303.5  # K34: was 1736934.2
$ws.Cells.Item(34, 12).Value = 1084.6666  # L34: was 1070.2222
$ws.Cells.Item(34, 13).Value = -1684101.5  # M34: was -1736732.2
$ws.Cells.Item(34, 14).Value = -1488.6666  # N34: was -1474.2222

$ws.Cells.Item(58, 8).Value = 32258758  # H58: was 31250690
$ws.Cells.Item(58, 9).Value = 45455096  # I58: was 37037656
$ws.Cells.Item(58, 10).Value = 1044.8889  # J58: was 1060.8
$ws.Cells.Item(58, 11).Value = 45455096  # K58: was 37037656
$ws.Cells.Item(58, 12).Value = 1044.8889  # L58: was 1060.8
$ws.Cells.Item(58, 13).Value = -45454893  # M58: was -37037453
$ws.Cells.Item(58, 14).Value = -1450.8889  # N58: was -1466.8

$ws.Cells.Item(122, 8).Value = 13889676  # H122: was 17858028
$ws.Cells.Item(122, 9).Value = 25000394  # I122: was 35714696
$ws.Cells.Item(122, 10).Value = 1279.25  # J122: was 1360.5714
$ws.Cells.Item(122, 11).Value = 75001182  # K122: was 107144088
$ws.Cells.Item(122, 12).Value = 3837.75  # L122: was 4081.7142
$ws.Cells.Item(122, 13).Value = -74998732  # M122: was -107141638
$ws.Cells.Item(122, 14).Value = -8737.75  # N122: was -8981.7142

$ws.Cells.Item(132, 8).Value = 10102708  # H132: was 33335806
$ws.Cells.Item(132, 9).Value = 1172.1111  # I132: was 1815.8572
$ws.Cells.Item(132, 10).Value = 22224550  # J132: was 111115110
$ws.Cells.Item(132, 11).Value = 3516.3333  # K132: was 5447.571599999999
$ws.Cells.Item(132, 12).Value = 66673650  # L132: was 333345330
$ws.Cells.Item(132, 13).Value = -986.3333000000002  # M132: was -2917.571599999999
$ws.Cells.Item(132, 14).Value = -66678710  # N132: was -333350390

$ws.Cells.Item(134, 8).Value = 1072.6923  # H134: was 1118.2858
$ws.Cells.Item(134, 9).Value = 981  # I134: was 989.9
$ws.Cells.Item(134, 10).Value = 1219.4  # J134: was 1439.25
$ws.Cells.Item(134, 11).Value = 2943  # K134: was 2969.7
$ws.Cells.Item(134, 12).Value = 3658.2  # L134: was 4317.75
$ws.Cells.Item(134, 13).Value = -408  # M134: was -434.6999999999998
$ws.Cells.Item(134, 14).Value = -8728.200000000001  # N134: was -9387.75

$ws.Cells.Item(136, 8).Value = 32258758  # H136: was 31250690
$ws.Cells.Item(136, 9).Value = 45455096  # I136: was 37037656
$ws.Cells.Item(136, 10).Value = 1044.8889  # J136: was 1060.8
$ws.Cells.Item(136, 11).Value = 136365288  # K136: was 111112968
$ws.Cells.Item(136, 12).Value = 3134.6667  # L136: was 3182.4
$ws.Cells.Item(136, 13).Value = -136362738  # M136: was -111110418
$ws.Cells.Item(136, 14).Value = -8234.6667  # N136: was -8282.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 58829880  # H107: was 71118730
$ws.Cells.Item(107, 10).Value = 100010630  # J107: was 133347400
$ws.Cells.Item(107, 12).Value = 300031890  # L107: was 400042200
$ws.Cells.Item(107, 14).Value = -300035730  # N107: was -400046040

$ws.Cells.Item(123, 8).Value = 1623.625  # H123: was 2215
$ws.Cells.Item(123, 9).Value = 1623.625  # I123: was 2215
$ws.Cells.Item(123, 11).Value = 4870.875  # K123: was 6645
$ws.Cells.Item(123, 13).Value = -2420.875  # M123: was -4195

$ws.Cells.Item(131, 8).Value = 784.12  # H131: was 780.75555
$ws.Cells.Item(131, 10).Value = 818.67413  # J131: was 819.2152
$ws.Cells.Item(131, 12).Value = 2456.02239  # L131: was 2457.6456
$ws.Cells.Item(131, 14).Value = -12536.02239  # N131: was -12537.6456

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1125  # H113: was 2671
$ws.Cells.Item(113, 9).Value = 1125  # I113: was 1000
$ws.Cells.Item(113, 10).Value = 0  # J113: was 6013
$ws.Cells.Item(113, 11).Value = 1125  # K113: was 1000
$ws.Cells.Item(113, 12).Value = 0  # L113: was 6013
$ws.Cells.Item(113, 13).ClearContents()  # M113: was 1170 -> empty
$ws.Cells.Item(113, 14).Value = 1045  # N113: was -10353

$ws.Cells.Item(120, 8).Value = 29316.5  # H120: was 30158.5
$ws.Cells.Item(120, 10).Value = 29316.5  # J120: was 30158.5
$ws.Cells.Item(120, 12).Value = 29316.5  # L120: was 30158.5
$ws.Cells.Item(120, 14).Value = -38992.5  # N120: was -39834.5

$ws.Cells.Item(132, 8).Value = 2163.037  # H132: was 2061.1936
$ws.Cells.Item(132, 9).Value = 1821.6  # I132: was 1641.625
$ws.Cells.Item(132, 10).Value = 2589.8333  # J132: was 3499.7144
$ws.Cells.Item(132, 11).Value = 5464.799999999999  # K132: was 4924.875
$ws.Cells.Item(132, 12).Value = 7769.499899999999  # L132: was 10499.1432
$ws.Cells.Item(132, 13).Value = -2934.799999999999  # M132: was -2394.875
$ws.Cells.Item(132, 14).Value = -12829.4999  # N132: was -15559.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 1000  # H30: was 700
$ws.Cells.Item(30, 9).Value = 1000  # I30: was 700
$ws.Cells.Item(30, 11).Value = 1000  # K30: was 700
$ws.Cells.Item(30, 13).Value = -892  # M30: was -592

$ws.Cells.Item(39, 8).Value = 4693.3335  # H39: was 60000
$ws.Cells.Item(39, 9).Value = 3000  # I39: was 0
$ws.Cells.Item(39, 10).Value = 8080  # J39: was 60000
$ws.Cells.Item(39, 11).Value = 3000  # K39: was 0
$ws.Cells.Item(39, 12).Value = 8080  # L39: was 60000
$ws.Cells.Item(39, 13).Value = -2540  # M39: was empty
$ws.Cells.Item(39, 14).Value = -9000  # N39: was -60920

$ws.Cells.Item(58, 8).Value = 5900  # H58: was 0
$ws.Cells.Item(58, 9).Value = 2000  # I58: was 0
$ws.Cells.Item(58, 10).Value = 9800  # J58: was 0
$ws.Cells.Item(58, 11).Value = 2000  # K58: was 0
$ws.Cells.Item(58, 12).Value = 9800  # L58: was 0
$ws.Cells.Item(58, 13).Value = -1740  # M58: was empty
$ws.Cells.Item(58, 14).Value = -10320  # N58: was empty

$ws.Cells.Item(100, 8).Value = 2528.2856  # H100: was 2471.1428
$ws.Cells.Item(100, 9).Value = 3300.6  # I100: was 2375.375
$ws.Cells.Item(100, 10).Value = 2399.5667  # J100: was 2499.5186
$ws.Cells.Item(100, 11).Value = 3300.6  # K100: was 2375.375
$ws.Cells.Item(100, 12).Value = 2399.5667  # L100: was 2499.5186
$ws.Cells.Item(100, 13).Value = -2759.6  # M100: was -1834.375
$ws.Cells.Item(100, 14).Value = -3481.5667  # N100: was -3581.5186

$ws.Cells.Item(121, 8).Value = 0  # H121: was 50460
$ws.Cells.Item(121, 10).Value = 0  # J121: was 50460
$ws.Cells.Item(121, 12).ClearContents()  # L121: was 50460 -> empty
$ws.Cells.Item(121, 14).Value = 0  # N121: was -53954

$ws.Cells.Item(132, 8).Value = 32661248  # H132: was 28578664
$ws.Cells.Item(132, 9).Value = 57144764  # I132: was 40817932
$ws.Cells.Item(132, 10).Value = 16560.133  # J132: was 20375.25
$ws.Cells.Item(132, 11).Value = 171434292  # K132: was 122453796
$ws.Cells.Item(132, 12).Value = 49680.399  # L132: was 61125.75
$ws.Cells.Item(132, 13).Value = -171431762  # M132: was -122451266
$ws.Cells.Item(132, 14).Value = -54740.399  # N132: was -66185.75

$ws.Cells.Item(136, 8).Value = 102042870  # H136: was 102042910
$ws.Cells.Item(136, 9).Value = 71431070  # I136: was 60152640
$ws.Cells.Item(136, 10).Value = 200000620  # J136: was 500000500
$ws.Cells.Item(136, 11).Value = 214293210  # K136: was 180457920
$ws.Cells.Item(136, 12).Value = 600001860  # L136: was 1500001500
$ws.Cells.Item(136, 13).Value = -214290660  # M136: was -180455370
$ws.Cells.Item(136, 14).Value = -600006960  # N136: was -1500006600

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(13, 8).Value = 0  # H13: was 200
$ws.Cells.Item(13, 9).Value = 0  # I13: was 200
$ws.Cells.Item(13, 11).Value = 0  # K13: was 200
$ws.Cells.Item(13, 13).ClearContents()  # M13: was -60 -> empty

$ws.Cells.Item(121, 8).Value = 30420  # H121: was 0
$ws.Cells.Item(121, 10).Value = 30420  # J121: was 0
$ws.Cells.Item(121, 12).Value = 30420  # L121: was 0
$ws.Cells.Item(121, 14).Value = -33914  # N121: was empty

$ws.Cells.Item(132, 8).Value = 30398.777  # H132: was 33060.605
$ws.Cells.Item(132, 9).Value = 62487.65  # I132: was 62470.707
$ws.Cells.Item(132, 10).Value = 1687.6842  # J132: was 1812.375
$ws.Cells.Item(132, 11).Value = 187462.95  # K132: was 187412.121
$ws.Cells.Item(132, 12).Value = 5063.0526  # L132: was 5437.125
$ws.Cells.Item(132, 13).Value = -184932.95  # M132: was -184882.121
$ws.Cells.Item(132, 14).Value = -10123.0526  # N132: was -10497.125

$ws.Cells.Item(136, 8).Value = 8334865  # H136: was 9260920
$ws.Cells.Item(136, 9).Value = 14286549  # I136: was 16667623
$ws.Cells.Item(136, 10).Value = 2507  # J136: was 2540.6667
$ws.Cells.Item(136, 11).Value = 42859647  # K136: was 50002869
$ws.Cells.Item(136, 12).Value = 7521  # L136: was 7622.000100000001
$ws.Cells.Item(136, 13).Value = -42857097  # M136: was -50000319
$ws.Cells.Item(136, 14).Value = -12621  # N136: was -12722.0001
